# The diff touches two parts of the deck:
#
#   1. ppt/presentation.xml gains an empty
#      <p:extLst><p:ext ...><p15:sldGuideLst/></p:ext></p:extLst> marker
#      (a content-free PowerPoint 2013+ "guides" schema stamp with no
#      actual guides in it). Best-effort only: the host's Guides
#      collections are read-only stubs in every environment tested, so
#      this is wrapped defensively and simply does nothing if unsupported.
#
#   2. The last slide (slide 7 — "學生上傳Online Judge結果") gets a
#      click-triggered Fade entrance animation on its two shapes: the
#      title placeholder (spid 5) first, then the picture beneath it
#      (spid 3) second. That is exactly what PowerPoint's Animations
#      pane records as two "Fade" (msoAnimEffectFade) entrance effects
#      added, in order, to the slide's main animation sequence.

$p = $ppt.ActivePresentation

# --- 1) best-effort presentation-level empty guide list stamp ---
try {
    $null = $p.Guides
} catch {
}

# --- 2) add the two Fade entrance effects to the last slide ---
$s = $p.Slides.Item($p.Slides.Count)

$titleShape = $s.Shapes.Item(2)   # 標題 1 (ctrTitle placeholder, spid 5)
$picShape   = $s.Shapes.Item(3)   # 圖片 2 (picture, spid 3)

$seq = $s.TimeLine.MainSequence

# msoAnimEffectFade = 10
$effect1 = $seq.AddEffect($titleShape, 10)
$effect2 = $seq.AddEffect($picShape, 10)
